$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25; Excel shifts rows 25..119 down to 26..120
# and the sheet dimension grows from A1:R119 to A1:R120 automatically.
$ws.Rows("25").Insert()

# Populate the newly inserted row 25 with the new weekly price record
# (same market/category metadata as the surrounding rows, new date & prices).
$ws.Range("A25").Value2 = 3
$ws.Range("B25").Value2 = "Femacal de La Calera"
$ws.Range("C25").Value2 = "Coquimbo"
$ws.Range("D25").Value2 = 44575
$ws.Range("E25").Value2 = 5
$ws.Range("F25").Value2 = 100112052
$ws.Range("G25").Value2 = "Albahaca"
$ws.Range("H25").Value2 = "Sin especificar"
$ws.Range("I25").Value2 = "Primera"
$ws.Range("J25").Value2 = 110
$ws.Range("K25").Value2 = 4500
$ws.Range("L25").Value2 = 5000
$ws.Range("M25").Value2 = 4727
$ws.Range("N25").Value2 = '$/docena de matas'
$ws.Range("O25").Value2 = "Provincia de Quillota"
$ws.Range("P25").Value2 = 788
$ws.Range("Q25").Value2 = 6
$ws.Range("R25").Value2 = "Hortaliza"

# Keep the D column's date number format consistent with the rest of the column.
$ws.Range("D25").NumberFormat = $ws.Range("D26").NumberFormat
